$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-8 with the new consolidated token text values
$ws.Range("A2").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A3").Value = "('Bird', ['Token Creature — Bird', 'Flying', '3/3'])"
$ws.Range("A4").Value = "('Pentavite', ['Token Artifact Creature — Pentavite', 'Flying', '1/1'])"
$ws.Range("A5").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"
$ws.Range("A6").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A7").Value = "('Wurm', ['Token Creature — Wurm', '6/6'])"
$ws.Range("A8").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"

# Remove the now-obsolete rows 9-24 (their data was folded into rows 2-8 above)
$ws.Range("A9:A24").ClearContents()
